$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing data rows (old "Neutrophils" and "Resolving-Mac"
# sending-cluster rows) -- the refreshed TPM run only produced 3 data rows.
$ws.Rows("5:6").Delete() | Out-Null

# Row 2 (was "ECs" sending cluster) -> now "FAPs", with refreshed TPM values.
$ws.Range("A2").Value = "FAPs"
$ws.Range("G2").Value = 0.134031
$ws.Range("H2").Value = 0.402093
$ws.Range("I2").Value = 0.1466544264074474
$ws.Range("J2").Value = 0.1466544264074474
$ws.Range("Q2").Value = 0.03511880261999999
$ws.Range("R2").Value = 0.31606922358
$ws.Range("S2").Value = 0.1466544264074474
$ws.Range("T2").Value = 0.1466544264074474

# Row 3 (was "FAPs" sending cluster) -> now "MuSCs", with refreshed TPM values.
$ws.Range("A3").Value = "MuSCs"
$ws.Range("G3").Value = 0.080957
$ws.Range("H3").Value = 0.242871
$ws.Range("I3").Value = 0.08858176391034703
$ws.Range("J3").Value = 0.08858176391034703
$ws.Range("Q3").Value = 0.02121235314
$ws.Range("R3").Value = 0.19091117826
$ws.Range("S3").Value = 0.08858176391034703
$ws.Range("T3").Value = 0.08858176391034703

# Row 4 (was "MuSCs" sending cluster) -> now "Neutrophils", with refreshed TPM values.
$ws.Range("A4").Value = "Neutrophils"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.6989360000000001
$ws.Range("H4").Value = 2.096808
$ws.Range("I4").Value = 0.7647638096822056
$ws.Range("J4").Value = 0.7647638096822056
$ws.Range("Q4").Value = 0.18313521072
$ws.Range("R4").Value = 1.64821689648
$ws.Range("S4").Value = 0.7647638096822056
$ws.Range("T4").Value = 0.7647638096822056
